$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 33
$ws.Range("I2").Value = 138
$ws.Range("J2").Value = 560
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 158
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 116
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 70
$ws.Range("T2").Value = 95
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 855
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 835
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 14
$ws.Range("AA2").Value = 3
